$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text representation
# (values like "3.556.03" / "1.00" are display strings, not numbers),
# so force text formatting before writing the updated price strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.125.86"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.550.43"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "604.07"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "143.74"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "3.547.97"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").Value = "7.77"
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "4.152.84"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "0.0000205"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "29.92"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "3.529.96"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "66.175.62"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("E19").Value = "  +5.38%  "
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "14.64"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "429.61"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").Value = "79.79"
$ws.Range("D25").Value = "3.694.11"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").Value = "2.48"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "9.06"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").Value = "7.78"
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "3.547.36"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "25.33"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("D37").Value = "7.77"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").Value = "5.51"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").Value = "174.73"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("D41").Value = "0.0843"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").Value = "5.16"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "0.885"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "1.90"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "46.05"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "1.18"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("D48").Value = "24.62"
$ws.Range("E48").Value = "  -5.25%  "
$ws.Range("D49").Value = "2.38"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "7.09"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "22.88"
$ws.Range("E51").Value = "  +1.00%  "
